$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '60.914.30'
Set-TextValue 2 5 '  -3.64%  '
Set-TextValue 3 4 '2.909.18'
Set-TextValue 3 5 '  -4.67%  '
Set-TextValue 4 5 '  -0.13%  '
Set-TextValue 5 4 '586.56'
Set-TextValue 5 5 '  -1.91%  '
Set-TextValue 6 4 '146.26'
Set-TextValue 6 5 '  -4.24%  '
Set-TextValue 7 5 '  -0.01%  '
Set-TextValue 8 4 '0.502'
Set-TextValue 8 5 '  -3.57%  '
Set-TextValue 9 4 '2.908.52'
Set-TextValue 9 5 '  -4.62%  '
Set-TextValue 10 4 '6.72'
Set-TextValue 10 5 '  +5.42%  '
Set-TextValue 11 4 '0.144'
Set-TextValue 11 5 '  -6.20%  '
Set-TextValue 12 4 '0.447'
Set-TextValue 12 5 '  -3.33%  '
Set-TextValue 13 4 '0.0000224'
Set-TextValue 13 5 '  -4.92%  '
Set-TextValue 14 4 '33.39'
Set-TextValue 14 5 '  -5.15%  '
Set-TextValue 15 5 '  +0.51%  '
Set-TextValue 16 4 '3.387.08'
Set-TextValue 16 5 '  -4.76%  '
Set-TextValue 17 4 '60.799.33'
Set-TextValue 17 5 '  -3.87%  '
Set-TextValue 18 4 '6.75'
Set-TextValue 18 5 '  -4.38%  '
Set-TextValue 19 4 '2.904.70'
Set-TextValue 19 5 '  -5.00%  '
Set-TextValue 20 4 '426.43'
Set-TextValue 20 5 '  -6.28%  '
Set-TextValue 21 4 '13.55'
Set-TextValue 21 5 '  -5.38%  '
Set-TextValue 22 4 '0.671'
Set-TextValue 22 5 '  -3.62%  '
Set-TextValue 23 4 '7.11'
Set-TextValue 23 5 '  -5.63%  '
Set-TextValue 24 4 '79.97'
Set-TextValue 24 5 '  -3.57%  '
Set-TextValue 25 4 '11.00'
Set-TextValue 25 5 '  +1.07%  '
Set-TextValue 26 4 '2.21'
Set-TextValue 26 5 '  -3.15%  '
Set-TextValue 27 4 '11.83'
Set-TextValue 27 5 '  -3.67%  '
Set-TextValue 28 5 '  +0.03%  '
Set-TextValue 29 4 '0.999'
Set-TextValue 29 5 '  -0.20%  '
Set-TextValue 30 4 '7.24'
Set-TextValue 30 5 '  -3.04%  '
Set-TextValue 31 5 '  -4.29%  '
Set-TextValue 32 4 '2.17'
Set-TextValue 32 5 '  -0.23%  '
Set-TextValue 33 4 '26.37'
Set-TextValue 33 5 '  -5.10%  '
Set-TextValue 34 5 '  -4.77%  '
Set-TextValue 35 4 '0.0₃0840'
Set-TextValue 35 5 '  -2.17%  '
Set-TextValue 36 5 '  -3.09%  '
Set-TextValue 37 4 '5.62'
Set-TextValue 37 5 '  -5.14%  '
Set-TextValue 38 4 '49.43'
Set-TextValue 38 5 '  -2.37%  '
Set-TextValue 39 4 '2.96'
Set-TextValue 39 5 '  -5.80%  '
Set-TextValue 40 5 '  -5.10%  '
Set-TextValue 41 5 '  -0.40%  '
Set-TextValue 42 4 '8.66'
Set-TextValue 42 5 '  -5.21%  '
Set-TextValue 43 5 '  -1.71%  '
Set-TextValue 44 4 '41.50'
Set-TextValue 44 5 '  +1.38%  '
$ws.Cells.Item(45, 2).Value = 'Bittensor'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 45 4 '374.89'
Set-TextValue 45 5 '  -5.60%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 46 4 '0.0346'
Set-TextValue 46 5 '  -3.23%  '
Set-TextValue 47 4 '2.668.90'
Set-TextValue 47 5 '  -3.11%  '
Set-TextValue 48 4 '133.11'
Set-TextValue 48 5 '  -0.25%  '
Set-TextValue 50 4 '24.96'
Set-TextValue 50 5 '  +3.13%  '
Set-TextValue 51 5 '  -2.27%  '
